$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column D to fit the longer employee name ---
$ws.Columns("D").ColumnWidth = 36.08984375

# --- Insert 4 new rows after the existing data row (16) so the table grows
#     from 1 data row to 5 data rows; this pushes the old rows 21/22
#     (signature block) down to 25/26, matching the target layout. ---
$ws.Rows("17:20").Insert()

# --- Copy the formatting of the (still unique) data row down into the
#     three freshly inserted rows so they inherit the same borders/number
#     formats used throughout the table body. ---
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J19").PasteSpecial(-4122)
$ws.Range("B16:J16").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header block ---
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("E7").Value = "TRANSCASTAÑEDA S.A.S."
$ws.Range("B9").Value = "NIT"
$ws.Range("E9").Value = 8040078147
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 264214
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("C13").Value = 2
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 5

# --- Table header row ---
$ws.Range("B15").Value = "Tipo Doc Trabajador"
$ws.Range("C15").Value = "N° Doc Trabajador"
$ws.Range("D15").Value = "Nombre Trabajador"
$ws.Range("E15").Value = "Periodo Mora"
$ws.Range("F15").Value = "Valor Mora"
$ws.Range("G15").Value = "Salario Basico"
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# --- Data rows: new employee DANIEL ARTURO DE LA OSSA SAUCEDO (4 periods)
#     followed by the pre-existing DEMOSTENES VIDES ZAMBRANO row, now last. ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1051636136"
$ws.Range("D16").Value = "DANIEL ARTURO DE LA OSSA SAUCEDO"
$ws.Range("E16").Value = "2111"
$ws.Range("F16").Value = 58015
$ws.Range("G16").Value = 1450380

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1051636136"
$ws.Range("D17").Value = "DANIEL ARTURO DE LA OSSA SAUCEDO"
$ws.Range("E17").Value = "2110"
$ws.Range("F17").Value = 58015
$ws.Range("G17").Value = 1450380

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1051636136"
$ws.Range("D18").Value = "DANIEL ARTURO DE LA OSSA SAUCEDO"
$ws.Range("E18").Value = "2109"
$ws.Range("F18").Value = 58015
$ws.Range("G18").Value = 1450380

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1051636136"
$ws.Range("D19").Value = "DANIEL ARTURO DE LA OSSA SAUCEDO"
$ws.Range("E19").Value = "2108"
$ws.Range("F19").Value = 58015
$ws.Range("G19").Value = 1450380

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1104124941"
$ws.Range("D20").Value = "DEMOSTENES VIDES ZAMBRANO"
$ws.Range("E20").Value = "2410"
$ws.Range("F20").Value = 32154
$ws.Range("G20").Value = 2679480

# --- Signature block (shifted down to rows 25/26 by the inserted rows) ---
$ws.Range("B25").Value = "___________________________________"
$ws.Range("H25").Value = "___________________________________"
$ws.Range("B26").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H26").Value = "FIRMA DEL REPRESENTANTE LEGAL"
